# "Updated code for ICW and VOE"
#
# - Assets sheet: add a new named constant row (TotalAssetLimit)
# - Findings sheet: add a new finding row (ICW_AssetsOver52787) for the
#   "assets exceed limit" manual-review message

$wb = $excel.ActiveWorkbook

$wsAssets   = $wb.Worksheets.Item("Assets")
$wsFindings = $wb.Worksheets.Item("Findings")

# --- Assets sheet: new row 38, "TotalAssetLimit" constant (A=B, like the others) ---
$wsAssets.Select()
$wsAssets.Range("A38").Value = "TotalAssetLimit"
$wsAssets.Range("B38").Value = "TotalAssetLimit"
$wsAssets.Range("B41").Select()

# --- Findings sheet: insert a new row 12 for the ICW asset-limit finding ---
$wsFindings.Select()
$wsFindings.Rows.Item(12).Insert()
$wsFindings.Range("A12").Value = "ICW_AssetsOver52787"
$wsFindings.Range("B12").Value = "Assets exceed the asset limit, manual review is required. "

# Keep Findings as the active/visible sheet (matches original workbook state)
$wsFindings.Select()
$wsFindings.Range("B13").Select()
